$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure price/volume columns keep their textual formatting (some new values,
# e.g. "1.00", "382.59", look like plain numbers and would otherwise be
# auto-converted to numeric cells by Excel's type inference).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '51.267.91'
$ws.Range("D3").Value = '2.964.59'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '382.59'
$ws.Range("E5").Value = '  +0.62%  '
$ws.Range("D6").Value = '103.45'
$ws.Range("E6").Value = '  -1.99%  '
$ws.Range("D7").Value = '0.540'
$ws.Range("E7").Value = '  -0.59%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.591'
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("D10").Value = '36.73'
$ws.Range("E10").Value = '  -1.83%  '
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").Value = '0.0843'
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("D13").Value = '3.431.05'
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("D14").Value = '18.08'
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").Value = '7.45'
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").Value = '2.961.51'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").Value = '0.996'
$ws.Range("E17").Value = '  +3.86%  '
$ws.Range("D18").Value = '51.213.64'
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("D19").Value = '3.21'
$ws.Range("E19").Value = '  -7.09%  '
$ws.Range("E20").Value = '  -4.15%  '
$ws.Range("D21").Value = '12.62'
$ws.Range("E21").Value = '  -4.60%  '
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").Value = '68.59'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '262.59'
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").Value = '2.92'
$ws.Range("E25").Value = '  +3.80%  '
$ws.Range("D26").Value = '8.38'
$ws.Range("E26").Value = '  +12.46%  '
$ws.Range("D27").Value = '7.87'
$ws.Range("E27").Value = '  +5.32%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("E29").Value = '  +8.92%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").Value = '25.75'
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("D33").Value = '0.0457'
$ws.Range("E33").Value = '  +4.85%  '
$ws.Range("D34").Value = '34.11'
$ws.Range("E34").Value = '  -1.20%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = '50.62'
$ws.Range("E35").Value = '  -3.58%  '
$ws.Range("B36").Value = 'Toncoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D36").Value = '2.06'
$ws.Range("E36").Value = '  -0.68%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  -2.04%  '
$ws.Range("E39").Value = '  -3.03%  '
$ws.Range("E40").Value = '  -4.18%  '
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("E42").Value = '  -3.15%  '
$ws.Range("D43").Value = '121.62'
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("D44").Value = '21.44'
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("D46").Value = '0.272'
$ws.Range("E46").Value = '  -2.77%  '
$ws.Range("E47").Value = '  +2.69%  '
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").Value = '2.016.27'
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("D50").Value = '0.0350'
$ws.Range("E50").Value = '  +5.99%  '
$ws.Range("E51").Value = '  +13.75%  '
